# Update "want to go" attendance counts (column F) for a handful of events
# on both the "展览" sheet and the "全部类型" (all types) aggregate sheet,
# mirroring the regenerated site data (commit 456a3b4).

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

# Map of row -> new value for column F
$updates = @{
    9  = 1315
    12 = 912
    14 = 508
    16 = 239
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}

$wb.Save()
